$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '64.407.66'
$ws.Range("E2").Value = '  +1.58%  '
$ws.Range("D3").Value = '3.174.80'
$ws.Range("E3").Value = '  +2.58%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '594.08'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +2.10%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '148.49'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  +2.23%  '
$ws.Range("E7").Value = '  +0.05%  '
$ws.Range("D8").Value = '3.161.56'
$ws.Range("E8").Value = '  +2.41%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.534'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  +1.38%  '
$ws.Range("E10").Value = '  +1.67%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '5.92'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  +5.36%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.461'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  +1.30%  '
$ws.Range("E13").Value = '  +1.56%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '37.71'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  +1.90%  '
$ws.Range("D15").Value = '3.692.56'
$ws.Range("E15").Value = '  +2.23%  '
$ws.Range("E16").Value = '  +0.17%  '
$ws.Range("E17").Value = '  +3.58%  '
$ws.Range("D18").Value = '64.198.19'
$ws.Range("E18").Value = '  +1.46%  '
$ws.Range("D19").Value = '3.160.75'
$ws.Range("E19").Value = '  +2.26%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '471.49'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  +2.48%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '14.58'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  +2.85%  '
$ws.Range("E22").Value = '  +2.54%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '7.70'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  +4.21%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '2.44'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  +14.49%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '13.29'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  +3.48%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '81.64'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  +0.61%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '10.12'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  +12.37%  '
$ws.Range("E28").Value = '  +0.10%  '
$ws.Range("E29").Value = '  +2.71%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '2.24'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  +2.62%  '
$ws.Range("B31").Value = 'FirstDigitalUSD'
$ws.Range("C31").Value = 'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd'
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '1.00'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  +0.13%  '
$ws.Range("B32").Value = 'NEARProtocol'
$ws.Range("C32").Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '7.28'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  +5.54%  '
$ws.Range("B33").Value = 'EthereumClassic'
$ws.Range("C33").Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '28.43'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  +6.96%  '
$ws.Range("B34").Value = 'Hedera'
$ws.Range("C34").Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '0.117'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  +5.68%  '
$ws.Range("D35").Value = '0.0₃0862'
$ws.Range("E35").Value = '  +2.31%  '
$ws.Range("E36").Value = '  +3.56%  '
$ws.Range("E37").Value = '  +4.16%  '
$ws.Range("E38").Value = '  +1.41%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '3.31'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  -2.48%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '467.72'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  +8.47%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '51.43'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  +2.45%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '9.33'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  +7.64%  '
$ws.Range("E43").Value = '  +10.00%  '
$ws.Range("E44").Value = '  +2.86%  '
$ws.Range("D45").Value = '2.910.27'
$ws.Range("E45").Value = '  +1.22%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '39.86'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  +11.48%  '
$ws.Range("E47").Value = '  +0.46%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '133.35'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  +6.79%  '
$ws.Range("E50").Value = '  +5.72%  '
$ws.Range("E51").Value = '  +1.46%  '